$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "user1"
$ws.Range("C2").Value = 123
$ws.Range("D2").Value = "filme 1"

# Row 3
$ws.Range("B3").Value = "user2"
$ws.Range("C3").Value = 123
$ws.Range("D3").Value = "filme 2"

# Row 4
$ws.Range("B4").Value = "user3"
$ws.Range("C4").Value = 123
$ws.Range("D4").Value = "filme 3"
$ws.Range("E4").ClearContents()
